# Auto-generated edit script: apply Durandal_Profits leve-profit recalculation
# updates across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 492.1154
$ws.Range("I92").Value = 484.65
$ws.Range("K92").Value = 484.65
$ws.Range("M92").Value = 763.35
$ws.Range("H98").Value = 26319810.0
$ws.Range("I98").Value = 3800.4688
$ws.Range("J98").Value = 166671860.0
$ws.Range("K98").Value = 3800.4688
$ws.Range("L98").Value = 166671860.0
$ws.Range("M98").Value = -2302.4688
$ws.Range("N98").Value = -166674856.0
$ws.Range("H103").Value = 71429340.0
$ws.Range("I103").Value = 844.6667
$ws.Range("J103").Value = 125000730.0
$ws.Range("K103").Value = 2534.0001
$ws.Range("L103").Value = 375002190.0
$ws.Range("M103").Value = -1948.0001
$ws.Range("N103").Value = -375003362.0
$ws.Range("H122").Value = 26319810.0
$ws.Range("I122").Value = 3800.4688
$ws.Range("J122").Value = 166671860.0
$ws.Range("K122").Value = 11401.4064
$ws.Range("L122").Value = 500015580.0
$ws.Range("M122").Value = -8951.4064
$ws.Range("N122").Value = -500020480.0
$ws.Range("H127").Value = 457.3889
$ws.Range("I127").Value = 288.07144
$ws.Range("J127").Value = 1050.0
$ws.Range("K127").Value = 864.21432
$ws.Range("L127").Value = 3150.0
$ws.Range("M127").Value = 4095.78568
$ws.Range("N127").Value = -13070.0
$ws.Range("H129").Value = 963.05884
$ws.Range("I129").Value = 322.18182
$ws.Range("J129").Value = 1269.5652
$ws.Range("K129").Value = 966.54546
$ws.Range("L129").Value = 3808.6956
$ws.Range("M129").Value = 4033.45454
$ws.Range("N129").Value = -13808.6956
$ws.Range("H131").Value = 1625.0
$ws.Range("I131").Value = 782.1429
$ws.Range("J131").Value = 3100.0
$ws.Range("K131").Value = 2346.4287
$ws.Range("L131").Value = 9300.0
$ws.Range("M131").Value = 2693.5713
$ws.Range("N131").Value = -19380.0
$ws.Range("H132").Value = 3903.8857
$ws.Range("I132").Value = 4026.3333
$ws.Range("J132").Value = 3169.2
$ws.Range("K132").Value = 12078.9999
$ws.Range("L132").Value = 9507.599999999999
$ws.Range("M132").Value = -9548.999899999999
$ws.Range("N132").Value = -14567.6
$ws.Range("H138").Value = 4861.364
$ws.Range("I138").Value = 1749.129
$ws.Range("J138").Value = 6958.7393
$ws.Range("K138").Value = 5247.387
$ws.Range("L138").Value = 20876.2179
$ws.Range("M138").Value = -107.3869999999997
$ws.Range("N138").Value = -31156.2179
$ws.Range("H141").Value = 1003.4667
$ws.Range("I141").Value = 878.3461
$ws.Range("J141").Value = 1816.75
$ws.Range("K141").Value = 2635.0383
$ws.Range("L141").Value = 5450.25
$ws.Range("M141").Value = 2544.9617
$ws.Range("N141").Value = -15810.25

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3105.5715
$ws.Range("I45").Value = 3153.0
$ws.Range("J45").Value = 3042.3333
$ws.Range("K45").Value = 3153.0
$ws.Range("L45").Value = 3042.3333
$ws.Range("M45").Value = -2776.0
$ws.Range("N45").Value = -3796.3333
$ws.Range("H132").Value = 18538124.0
$ws.Range("I132").Value = 25641636.0
$ws.Range("J132").Value = 68990.266
$ws.Range("K132").Value = 76924908.0
$ws.Range("L132").Value = 206970.798
$ws.Range("M132").Value = -76922378.0
$ws.Range("N132").Value = -212030.798

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1644.1364
$ws.Range("I105").Value = 1601.6154
$ws.Range("J105").Value = 1705.5555
$ws.Range("K105").Value = 1601.6154
$ws.Range("L105").Value = 1705.5555
$ws.Range("M105").Value = 145.3846000000001
$ws.Range("N105").Value = -5199.5555

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 21557.48
$ws.Range("I132").Value = 1090.5135
$ws.Range("J132").Value = 79809.62
$ws.Range("K132").Value = 3271.5405
$ws.Range("L132").Value = 239428.86
$ws.Range("M132").Value = -741.5405000000001
$ws.Range("N132").Value = -244488.86

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 17858438.0
$ws.Range("J131").Value = 20001190.0
$ws.Range("L131").Value = 60003570.0
$ws.Range("N131").Value = -60013650.0
$ws.Range("H133").Value = 6138.593
$ws.Range("I133").Value = 1298.8889
$ws.Range("J133").Value = 7106.533
$ws.Range("K133").Value = 3896.6667
$ws.Range("L133").Value = 21319.599
$ws.Range("M133").Value = 1163.3333
$ws.Range("N133").Value = -31439.599

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 2001599.2
$ws.Range("J21").Value = 1999.0
$ws.Range("L21").Value = 1999.0
$ws.Range("N21").Value = -2345.0
$ws.Range("H30").Value = 2001599.2
$ws.Range("J30").Value = 1999.0
$ws.Range("L30").Value = 1999.0
$ws.Range("N30").Value = -2209.0
$ws.Range("H33").Value = 616538.06
$ws.Range("J33").Value = 616538.06
$ws.Range("L33").Value = 616538.06
$ws.Range("N33").Value = -617042.06
$ws.Range("H99").Value = 4710.2666
$ws.Range("I99").Value = 3761.0715
$ws.Range("J99").Value = 17999.0
$ws.Range("K99").Value = 3761.0715
$ws.Range("L99").Value = 17999.0
$ws.Range("M99").Value = -1515.0715
$ws.Range("N99").Value = -22491.0
$ws.Range("H102").Value = 1333.3667
$ws.Range("I102").Value = 1519.6666
$ws.Range("J102").Value = 1053.9166
$ws.Range("K102").Value = 1519.6666
$ws.Range("L102").Value = 1053.9166
$ws.Range("M102").Value = 102.3334
$ws.Range("N102").Value = -4297.9166
$ws.Range("H108").Value = 25995.0
$ws.Range("J108").Value = 25995.0
$ws.Range("L108").Value = 25995.0
$ws.Range("N108").Value = -33675.0

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1000.0
$ws.Range("I22").Value = 1000.0
$ws.Range("J22").Value = 1000.0
$ws.Range("K22").Value = 1000.0
$ws.Range("L22").Value = 1000.0
$ws.Range("M22").Value = -705.0
$ws.Range("N22").Value = -1590.0
$ws.Range("H27").Value = 1000.0
$ws.Range("I27").Value = 1000.0
$ws.Range("J27").Value = 1000.0
$ws.Range("K27").Value = 1000.0
$ws.Range("L27").Value = 1000.0
$ws.Range("M27").Value = -893.0
$ws.Range("N27").Value = -1214.0
$ws.Range("H40").Value = 1913.762
$ws.Range("I40").Value = 1781.7059
$ws.Range("J40").Value = 2475.0
$ws.Range("K40").Value = 1781.7059
$ws.Range("L40").Value = 2475.0
$ws.Range("M40").Value = -1645.7059
$ws.Range("N40").Value = -2747.0
$ws.Range("H93").Value = 1243.8108
$ws.Range("I93").Value = 1221.9584
$ws.Range("J93").Value = 1284.1538
$ws.Range("K93").Value = 1221.9584
$ws.Range("L93").Value = 1284.1538
$ws.Range("M93").Value = 26.04160000000002
$ws.Range("N93").Value = -3780.1538

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 4000.0
$ws.Range("I40").Value = 1750.0
$ws.Range("J40").Value = 5800.0
$ws.Range("K40").Value = 1750.0
$ws.Range("L40").Value = 5800.0
$ws.Range("M40").Value = -1601.0
$ws.Range("N40").Value = -6098.0
$ws.Range("H136").Value = 35221.69
$ws.Range("I136").Value = 48098.57
$ws.Range("J136").Value = 1419.875
$ws.Range("K136").Value = 144295.71
$ws.Range("L136").Value = 4259.625
$ws.Range("M136").Value = -141745.71
$ws.Range("N136").Value = -9359.625
